# Add tips for when the citizen button is not pressed.
#
# Appends five new key/value rows (TIPS_CITIZEN_1..5) to the bottom of the
# "StringLocalizations_BasicText" sheet (sheet4), copying the formatting of
# the last existing row (row 68), and makes that sheet the active tab /
# selection, matching the author's commit.

$wb  = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("StringLocalizations_BasicText")

# --- Append 5 new rows, copying formatting (style + wrap) from the last row ---
$lastRow = $ws4.Range("A68:E68")
$lastRow.Copy($ws4.Range("A69:E69"))
$lastRow.Copy($ws4.Range("A70:E70"))
$lastRow.Copy($ws4.Range("A71:E71"))
$lastRow.Copy($ws4.Range("A72:E72"))
$lastRow.Copy($ws4.Range("A73:E73"))

# Keys first (column A), so the new shared strings come out key-block then
# value-block, matching the source order.
$ws4.Cells.Item(69, 1).Value2 = "TIPS_CITIZEN_1"
$ws4.Cells.Item(70, 1).Value2 = "TIPS_CITIZEN_2"
$ws4.Cells.Item(71, 1).Value2 = "TIPS_CITIZEN_3"
$ws4.Cells.Item(72, 1).Value2 = "TIPS_CITIZEN_4"
$ws4.Cells.Item(73, 1).Value2 = "TIPS_CITIZEN_5"

# Then the English (en-gb) copy, column B.
$ws4.Cells.Item(69, 2).Value2 = "TIP:*n*You can save resources by asking citizens for more information!"
$ws4.Cells.Item(70, 2).Value2 = "TIP:*n*Citizens using the INSPEC2T app can help to identify key suspects"
$ws4.Cells.Item(71, 2).Value2 = "TIP:*n*Citizens may be able to help with cases with information or evidence, make the most of the option"
$ws4.Cells.Item(72, 2).Value2 = "TIP:*n*Maybe next time ask citizens for information, it could cut down the required officers"
$ws4.Cells.Item(73, 2).Value2 = "TIP:*n*Maybe next time ask for citizen help from the INSPEC2T app, it could reduce the number of turns required"

# --- Update the active sheet / selection so it matches the saved workbook view ---
$ws4.Activate()
$ws4.Range("B73").Select()
